$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a single row at row 9. This pushes the old rows 9,10,11 (the section
# labels) down to 10,11,12 and leaves rows 7 and 8 free (they did not exist
# in the original sparse sheet) to be populated with two new scenario rows.
$ws.Rows.Item(9).Insert()

# copy row 6's number formats / fills down into the two new rows 7 and 8
$ws.Range("A6:R6").Copy()
$ws.Range("A7:R8").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Row 6 input changes (re-run of an existing scenario with new inputs) ---
$ws.Range("A6").Value = 0.013
$ws.Range("N6").Value = 0.55
$ws.Range("S6").Formula = "=P6*1"
$ws.Range("S6").Style = "Normal"

# --- New row 7 ---
$ws.Range("A7").Value = 0.012
$ws.Range("B7").Value = 0.01
$ws.Range("C7").Formula = "=A7-B7"
$ws.Range("D7").Formula = "=A7/B7-1"
$ws.Range("E7").Value = 0.05
$ws.Range("F7").Value = 0.8
$ws.Range("G7").Formula = "=(A7*P7+B7*Q7)/O7"
$ws.Range("H7").Formula = "=_xlfn.NORM.S.INV(1-E7)"
$ws.Range("I7").Formula = "=_xlfn.NORM.S.INV(F7)"
$ws.Range("J7").Formula = "=SQRT(G7*(1-G7)*(1/Q7+1/P7))"
$ws.Range("K7").Formula = "=(A7-B7)/J7"
$ws.Range("L7").Formula = "=H7+I7"
$ws.Range("M7").Formula = "=(_xlfn.NORM.S.DIST(K7,TRUE))"
$ws.Range("N7").Value = 0.8
$ws.Range("O7").Value = 70000
$ws.Range("P7").Formula = "=O7*N7"
$ws.Range("Q7").Formula = "=O7-P7"
$ws.Range("R7").Formula = "=IF(K7>L7,""Yes"",""No"")"
$ws.Range("T7").Formula = "=10*24*P7*C7"
$ws.Range("T7").Style = "Normal"

# --- New row 8 ---
$ws.Range("A8").Value = 0.012
$ws.Range("B8").Value = 0.01
$ws.Range("C8").Formula = "=A8-B8"
$ws.Range("D8").Formula = "=A8/B8-1"
$ws.Range("E8").Value = 0.05
$ws.Range("F8").Value = 0.8
$ws.Range("G8").Formula = "=(A8*P8+B8*Q8)/O8"
$ws.Range("H8").Formula = "=_xlfn.NORM.S.INV(1-E8)"
$ws.Range("I8").Formula = "=_xlfn.NORM.S.INV(F8)"
$ws.Range("J8").Formula = "=SQRT(G8*(1-G8)*(1/Q8+1/P8))"
$ws.Range("K8").Formula = "=(A8-B8)/J8"
$ws.Range("L8").Formula = "=H8+I8"
$ws.Range("M8").Formula = "=(_xlfn.NORM.S.DIST(K8,TRUE))"
$ws.Range("N8").Value = 0.9
$ws.Range("O8").Value = 70000
$ws.Range("P8").Formula = "=O8*N8"
$ws.Range("Q8").Formula = "=O8-P8"
$ws.Range("R8").Formula = "=IF(K8>L8,""Yes"",""No"")"

# --- view / selection updates ---
$ws.Application.ActiveWindow.Zoom = 115
$ws.Range("U10").Select()
